$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.3048080303191223
$ws.Range("C2").Value = 10.29869402782916
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 20.06922418068126

$ws.Range("B3").Value = 0.01514828764759746
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 3.900430680208489
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 6.080152761294947

$ws.Range("B4").Value = 1.459612070389937
$ws.Range("C4").Value = 1.667794583268128
$ws.Range("D4").Value = 0.1575252929769615
$ws.Range("E4").Value = 0.496779210170732
$ws.Range("G4").Value = 3.781711156805759

$ws.Range("B5").Value = 0.04763786555579896
$ws.Range("C5").Value = 0.04240448674262143
$ws.Range("D5").Value = 3.900430680208489
$ws.Range("E5").Value = 616238.5361209477
$ws.Range("G5").Value = 616242.5265939801

$ws.Range("B6").Value = 0.6753301551942219
$ws.Range("C6").Value = 0.3127903958511391
$ws.Range("D6").Value = 0.8054896365839992
$ws.Range("E6").Value = 8.660232485948974
$ws.Range("G6").Value = 10.45384267357833

$ws.Range("B7").Value = 0.3048080303191223
$ws.Range("C7").Value = 0.3127903958511391
$ws.Range("D7").Value = 0.8054896365839992
$ws.Range("E7").Value = 8.660232485948974
$ws.Range("G7").Value = 10.08332054870323

$ws.Range("B8").Value = 3.230985683306322
$ws.Range("C8").Value = 1.667794583268128
$ws.Range("D8").Value = 0.8054896365839992
$ws.Range("E8").Value = 0.496779210170732
$ws.Range("G8").Value = 6.201049113329182

$ws.Range("B9").Value = 3.230985683306322
$ws.Range("C9").Value = 1.667794583268128
$ws.Range("D9").Value = 0.8054896365839992
$ws.Range("E9").Value = 8.660232485948974
$ws.Range("G9").Value = 14.36450238910742

